# Reorder columns D (codeforiati:group-name), E (codeforiati:category-name),
# F (codeforiati:category-code) so that the category-code column moves to
# position D (right after "status"), and group-name / category-name shift
# one column to the right into E / F. Column G (codeforiati:group-code) is
# left untouched.
#
# Before: code | name | status | group-name | category-name | category-code | group-code
# After : code | name | status | category-code | group-name | category-name | group-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
$firstRow = $usedRange.Row

# The category-code values ("111", "112", ...) look like numbers. Force the
# destination column to be formatted as text *before* writing so Excel keeps
# storing them as text/shared-strings instead of silently converting them to
# numeric cells.
$ws.Range("D$firstRow", "D$lastRow").NumberFormat = "@"
$ws.Range("E$firstRow", "E$lastRow").NumberFormat = "@"
$ws.Range("F$firstRow", "F$lastRow").NumberFormat = "@"

for ($rowNum = $firstRow; $rowNum -le $lastRow; $rowNum++) {
    $oldGroupName    = $ws.Cells.Item($rowNum, 4).Value()
    $oldCategoryName = $ws.Cells.Item($rowNum, 5).Value()
    $oldCategoryCode = $ws.Cells.Item($rowNum, 6).Value()

    $ws.Cells.Item($rowNum, 4).Value = $oldCategoryCode   # D: codeforiati:category-code
    $ws.Cells.Item($rowNum, 5).Value = $oldGroupName      # E: codeforiati:group-name
    $ws.Cells.Item($rowNum, 6).Value = $oldCategoryName   # F: codeforiati:category-name
}
